# Update marksheet totals: correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Total row "Right" count: 3 -> 5
$ws.Range("B11").Value = 5

# Total row "Total" marks: 48 -> 80
$ws.Range("B12").Value = 80

# Corr/total marks display string: "44/84" -> "80/140"
$ws.Range("E12").Value = "80/140"
